# Generate Report for Handback
#
# Updates the localization-status report:
#   - Overview sheet: "Ready for handoff" status becomes
#     "Handed back: in sync with en-US" for both zh-cn / de-de columns.
#   - zh-cn / de-de sheets: refresh the "Latest Handback DateTime" stamps,
#     and clear the (now stale) "handback not latest" warning message.
#   - Widen the Status / Latest Handback DateTime columns and narrow the
#     Error Detail column to fit the new content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "Handed back: in sync with en-US"
$overview.Range("E1:F1").ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-08-12 20:43:14"
$zhcn.Range("K3").Value = "2016-08-12 20:43:14"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""
$zhcn.Range("C1").ColumnWidth = 29.166666666666664
$zhcn.Range("P1").ColumnWidth = 12.833333333333332

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 20:43:24"
$dede.Range("P2").Value = ""
$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 20:43:24"
$dede.Range("P3").Value = ""
$dede.Range("C1").ColumnWidth = 29.166666666666664
$dede.Range("P1").ColumnWidth = 12.833333333333332
